$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously-empty row 96 with a new employee record (Test Karna) ---
# Copy formatting from row 95 so the date-formatted cells (C/F) reuse the existing style.
$ws.Range("A95:F95").Copy()
$ws.Range("A96:F96").PasteSpecial(-4122)

$ws.Cells.Item(96, 1).Value = 201
$ws.Cells.Item(96, 2).Value = "Test Karna"
$ws.Cells.Item(96, 3).Value = 44753
$ws.Cells.Item(96, 4).Value = "Software Engineer"
$ws.Cells.Item(96, 5).Value = ".Net"
$ws.Cells.Item(96, 6).Value = 35800

# --- Clear the footer text ("Generated on .../This report is generated by...") ---
# while keeping the cell styles (s=6 / s=7) intact.
$ws.Cells.Item(97, 1).Value = ""
$ws.Cells.Item(97, 3).Value = ""

# --- Update the view/selection to match the edited workbook state ---
$excel.Goto($ws.Range("A91"), $true)
$ws.Range("D101").Select()
